# Update the workbook per the commit:
#   refactor: migrate supporting CAF as a pacakge
#
# The chemical-formula list on Sheet1 is edited:
#   - the "Th2Os" entry is renamed to "ThOs"
#   - the "Sn5Co2" entry is renamed/replaced with "YNdThSi2"
#   - rows A4/A5 swap which of those two values they display
#   - the active selection moves from A2:A4 to A6:XFD7 (whole rows 6-7)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the shared-string values themselves (affects every cell using them).
$ws.Cells.Replace("Th2Os", "ThOs", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
$ws.Cells.Replace("Sn5Co2", "YNdThSi2", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)

# Swap the contents of A4 and A5.
$a4 = $ws.Range("A4").Value2
$a5 = $ws.Range("A5").Value2
$ws.Range("A4").Value = $a5
$ws.Range("A5").Value = $a4

# Move the active selection to A6:XFD7 (i.e. rows 6-7 selected), active cell A6.
$ws.Rows("6:7").Select()
